$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: A(Sending cluster), D(Target cluster), E..T
$data = @(
  @("ECs","FAPs",3,1,3.075165666666667,9.225497000000001,0.02641273658732285,0.02641273658732285,3,1,1.183046666666667,3.54914,0.6222589862820888,0.6222589862820888,3.638064491397778,32.74258042258001,0.01643556269376335,0.01643556269376335),
  @("ECs","MuSCs",3,1,3.075165666666667,9.225497000000001,0.02641273658732285,0.02641273658732285,3,1,0.718166,2.154498,0.3777410137179113,0.3777410137179112,2.208479426167334,19.876314835506,0.009977173893559496,0.009977173893559496),
  @("FAPs","FAPs",3,1,63.967809,191.903427,0.549422396165273,0.5494223961652731,3,1,1.183046666666667,3.54914,0.6222589862820888,0.6222589862820888,75.67690321141998,681.09212890278,0.341883023278479,0.341883023278479),
  @("FAPs","MuSCs",3,1,63.967809,191.903427,0.549422396165273,0.5494223961652731,3,1,0.718166,2.154498,0.3777410137179113,0.3777410137179112,45.939505518294,413.455549664646,0.2075393728867941,0.2075393728867941),
  @("MuSCs","FAPs",3,1,49.38440333333333,148.15321,0.424164867247404,0.4241648672474041,3,1,1.183046666666667,3.54914,0.6222589862820888,0.6222589862820888,58.42405374882222,525.8164837394,0.2639404003098464,0.2639404003098464),
  @("MuSCs","MuSCs",3,1,49.38440333333333,148.15321,0.424164867247404,0.4241648672474041,3,1,0.718166,2.154498,0.3777410137179113,0.3777410137179112,35.46619940428666,319.19579463858,0.1602244669375577,0.1602244669375577)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]

  $ws.Cells.Item($r, 1).Value = $row[0]   # A Sending cluster
  $ws.Cells.Item($r, 2).Value = "Thbs1"   # B Ligand symbol
  $ws.Cells.Item($r, 3).Value = "Tnfrsf11b" # C Receptor symbol
  $ws.Cells.Item($r, 4).Value = $row[1]   # D Target cluster

  for ($c = 5; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 3]
  }
}

# Remove now-unused rows 8, 9, 10 (previously the sending cluster "MuSCs" paired with itself etc.)
$ws.Range("A8:T10").Delete()
